# New weekly Brócoli price record for "Macroferia Regional de Talca" was
# added to the logged series. In the source table this new observation is
# inserted as row 469, pushing every following row (old 469..599) down by
# one position (to 470..600), which is why dimension grows from R599 to R600.
#
# Insert a blank row at row 469 (shifts rows 469-599 -> 470-600, carrying
# their existing values/formatting down with them) and then populate the
# newly-inserted row with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("469:469").Insert()

$ws.Cells.Item(469, 1).Value = 5
$ws.Cells.Item(469, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(469, 3).Value = "Maule"
$ws.Cells.Item(469, 4).Value = 45135
$ws.Cells.Item(469, 5).Value = 7
$ws.Cells.Item(469, 6).Value = 100112023
$ws.Cells.Item(469, 7).Value = "Brócoli"
$ws.Cells.Item(469, 8).Value = "Sin especificar"
$ws.Cells.Item(469, 9).Value = "Primera"
$ws.Cells.Item(469, 10).Value = 5000
$ws.Cells.Item(469, 11).Value = 500
$ws.Cells.Item(469, 12).Value = 500
$ws.Cells.Item(469, 13).Value = 500
$ws.Cells.Item(469, 14).Value = "$/unidad"
$ws.Cells.Item(469, 15).Value = "Región del Maule"
$ws.Cells.Item(469, 16).Value = 500
$ws.Cells.Item(469, 17).Value = 1
$ws.Cells.Item(469, 18).Value = "Hortaliza"
